# Append " a Tropy" to the end of the "Jazykové prostředky" Nadpis2 heading,
# turning it into "Jazykové prostředky a Tropy".
$d = $word.ActiveDocument

$heading = $d.Content
$found = $heading.Find.Execute("Jazykové prostředky", $true, $false, $false,
                                $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $heading.Collapse(0)          # wdCollapseEnd -> collapse the found range to its end
    $heading.InsertAfter(" a Tropy")
}
